# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" (exhibition) sheet, column F
$exhibitUpdates = @{
    5  = 7067
    7  = 916
    11 = 11
    12 = 58
    16 = 2825
    17 = 139
    23 = 120
    25 = 102
    26 = 149
    32 = 242
    33 = 360
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value for the "全部类型" (all types) sheet, column F
$allTypesUpdates = @{
    9  = 7067
    11 = 916
    15 = 11
    16 = 58
    21 = 2825
    22 = 139
    30 = 120
    32 = 102
    33 = 149
    39 = 242
    40 = 360
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
